$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear out the old K:P columns (no longer used) ---
$ws.Range("K1:P7").Clear()

# --- Set column widths for A:J ---
# NOTE: Excel's ColumnWidth setter pads the stored OOXML <col width> value by a
# constant +5/6 (0.8333333333333334) versus the character width that is set.
# Subtract that padding here so the saved width matches the target exactly.
$padding = 0.8333333333333334
$ws.Range("A1").EntireColumn.ColumnWidth = 44 - $padding
$ws.Range("B1").EntireColumn.ColumnWidth = 44 - $padding
$ws.Range("C1").EntireColumn.ColumnWidth = 43 - $padding
$ws.Range("D1").EntireColumn.ColumnWidth = 27 - $padding
$ws.Range("E1").EntireColumn.ColumnWidth = 26 - $padding
$ws.Range("F1").EntireColumn.ColumnWidth = 25 - $padding
$ws.Range("G1").EntireColumn.ColumnWidth = 11 - $padding
$ws.Range("H1").EntireColumn.ColumnWidth = 14 - $padding
$ws.Range("I1").EntireColumn.ColumnWidth = 13 - $padding
$ws.Range("J1").EntireColumn.ColumnWidth = 12 - $padding

# --- Header row ---
$ws.Range("A1").Value = "Filename"
$ws.Range("B1").Value = "PSNR Ground checker diff Reference checker"
$ws.Range("C1").Value = "PSNR Ground checker diff Enhanced checker"
$ws.Range("D1").Value = "MBE Ground diff Reference"
$ws.Range("E1").Value = "MBE Ground diff Enhanced"
$ws.Range("F1").Value = "MBE Ground diff Dehazed"
$ws.Range("G1").Value = "AG Ground"
$ws.Range("H1").Value = "AG Reference"
$ws.Range("I1").Value = "AG Enhanced"
$ws.Range("J1").Value = "AG Dehazed"

# --- Data rows (Filenames are unchanged) ---
$data = @{
    2 = @(8.57, 10.01, 6.07, 19.83, 20.21, 8.619999999999999, 4.77, 24.03, 5.54)
    3 = @(7.27, 8.85, 13.03, 17.27, 22.58, 8.619999999999999, 3.42, 30.44, 4.11)
    4 = @(8.77, 6.98, -3.96, 19.51, 18.49, 8.619999999999999, 7.92, 12.11, 8.539999999999999)
    5 = @(7.12, 5.99, 8.6, 17.57, 22.04, 8.619999999999999, 4.81, 16.73, 5.48)
    6 = @(11.48, 10.17, -4.13, 6.19, 16.14, 8.619999999999999, 7.73, 13.82, 8.24)
    7 = @(8.74, 10.09, 8.539999999999999, 4.46, 21, 8.619999999999999, 4.79, 16.88, 5.51)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = $i + 2  # B = 2
        $ws.Cells.Item($r, $col).Value = $vals[$i]
    }
}
